# Switch license from BY-NC to BY-SA
# Plus a couple of small copy-merges that accompanied the edit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("Blue Waters Petascale ... Semester Curriculum v1.0" title slide)
# The title was originally typed as two runs ("Blue Waters Petascale" and
# " Semester Curriculum v1.0") with identical formatting; merge them into a
# single run/string, matching a retype of that line.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

$found = $titleRange.Find("Blue Waters Petascale Semester Curriculum v1.0", 0)
$mergedTitle = $titleRange.Characters($found.Start, $found.Length)
$mergedTitle.Text = "Blue Waters Petascale Semester Curriculum v1.0"

# ---------------------------------------------------------------------------
# Slide 2 (license / colophon slide)
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$licenseShape = $slide2.Shapes.Item(1)
$licenseRange = $licenseShape.TextFrame.TextRange

# Nudge the placeholder's horizontal position by 2 EMU (754745 -> 754743).
$licenseShape.Left = 59.4286

# "CC BY-NC 4.0. To view a copy of this license, visit " ->
#   "CC " + "BY-SA " + "4.0. To view a copy of this license, visit "
$ccFound = $licenseRange.Find("CC BY-NC 4.0. To view a copy of this license, visit ", 0)
$ccStart = $ccFound.Start
$byNc = $licenseRange.Characters($ccStart + 3, 6)
$byNc.Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0" ->
#   "https://" + "creativecommons.org/licenses/by-sa/4.0"  (keeps the hyperlink)
$linkFound = $licenseRange.Find("https://creativecommons.org/licenses/by-nc/4.0", 0)
$linkStart = $linkFound.Start
$linkRest = $licenseRange.Characters($linkStart + 8, 38)
$linkRest.Text = "creativecommons.org/licenses/by-sa/4.0"

# ---------------------------------------------------------------------------
# Slide 27 (Acknowledgements) - merge "2) " + "The " into a single run.
# ---------------------------------------------------------------------------
$slide27 = $p.Slides.Item(27)
$ackShape = $slide27.Shapes.Item(2)
$ackRange = $ackShape.TextFrame.TextRange

$twoFound = $ackRange.Find("2) The ", 0)
$mergedTwo = $ackRange.Characters($twoFound.Start, $twoFound.Length)
$mergedTwo.Text = "2) The "
